$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '70.836.79'
$ws.Range("E2").Value = '  +1.57%  '
$ws.Range("D3").Value = '3.644.59'
$ws.Range("E3").Value = '  +3.48%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '604.38'
$ws.Range("E5").Value = '  -0.61%  '
$ws.Range("D6").Value = '197.82'
$ws.Range("E6").Value = '  +0.86%  '
$ws.Range("D7").Value = '0.628'
$ws.Range("E7").Value = '  -0.35%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").Value = '0.218'
$ws.Range("E9").Value = '  +8.55%  '
$ws.Range("D10").Value = '0.646'
$ws.Range("E10").Value = '  -0.42%  '
$ws.Range("D11").Value = '53.78'
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("E12").Value = '  +1.57%  '
$ws.Range("E13").Value = '  +0.59%  '
$ws.Range("D14").Value = '4.216.16'
$ws.Range("E14").Value = '  +3.27%  '
$ws.Range("D15").Value = '606.62'
$ws.Range("E15").Value = '  +1.65%  '
$ws.Range("D16").Value = '13.03'
$ws.Range("E16").Value = '  +1.59%  '
$ws.Range("D17").Value = '70.928.84'
$ws.Range("E17").Value = '  +1.44%  '
$ws.Range("D18").Value = '3.611.86'
$ws.Range("E18").Value = '  +2.22%  '
$ws.Range("D19").Value = '19.06'
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("E20").Value = '  +1.11%  '
$ws.Range("D21").Value = '1.00'
$ws.Range("E21").Value = '  +0.79%  '
$ws.Range("D22").Value = '18.26'
$ws.Range("E22").Value = '  +1.05%  '
$ws.Range("D23").Value = '5.30'
$ws.Range("E23").Value = '  -0.32%  '
$ws.Range("D24").Value = '103.97'
$ws.Range("E24").Value = '  +1.62%  '
$ws.Range("D25").Value = '4.63'
$ws.Range("E25").Value = '  -0.80%  '
$ws.Range("D26").Value = '3.00'
$ws.Range("E26").Value = '  -5.32%  '
$ws.Range("D27").Value = '10.63'
$ws.Range("E27").Value = '  -2.18%  '
$ws.Range("D28").Value = '9.75'
$ws.Range("E28").Value = '  +1.43%  '
$ws.Range("D29").Value = '33.90'
$ws.Range("E29").Value = '  +1.19%  '
$ws.Range("D30").Value = '4.73'
$ws.Range("E30").Value = '  +11.80%  '
$ws.Range("D31").Value = '7.24'
$ws.Range("E31").Value = '  +2.75%  '
$ws.Range("D32").Value = '12.31'
$ws.Range("E32").Value = '  -1.07%  '
$ws.Range("E33").Value = '  +0.66%  '
$ws.Range("D34").Value = '63.45'
$ws.Range("E34").Value = '  +0.48%  '
$ws.Range("D35").Value = '0.0₃0891'
$ws.Range("E35").Value = '  +4.93%  '
$ws.Range("D36").Value = '3.980.53'
$ws.Range("E36").Value = '  +7.19%  '
$ws.Range("E37").Value = '  +0.39%  '
$ws.Range("E38").Value = '  +0.35%  '
$ws.Range("D39").Value = '517.62'
$ws.Range("E39").Value = '  +5.39%  '
$ws.Range("D40").Value = '0.391'
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("D41").Value = '36.71'
$ws.Range("E41").Value = '  +0.44%  '
$ws.Range("D42").Value = '3.56'
$ws.Range("E42").Value = '  -2.36%  '
$ws.Range("E43").Value = '  +2.77%  '
$ws.Range("D44").Value = '0.0462'
$ws.Range("E44").Value = '  +1.62%  '
$ws.Range("D45").Value = '3.47'
$ws.Range("E45").Value = '  +5.68%  '
$ws.Range("D46").Value = '2.92'
$ws.Range("E46").Value = '  +4.00%  '
$ws.Range("D47").Value = '0.141'
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("D48").Value = '8.60'
$ws.Range("E48").Value = '  +1.53%  '
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("D50").Value = '0.000252'
$ws.Range("E50").Value = '  +2.53%  '
$ws.Range("D51").Value = '1.31'
$ws.Range("E51").Value = '  +0.91%  '
